# Apply the Sat Jun 10 2023 GitHub Actions "Updated cryptos list" refresh:
# new Price (D) / Volume(1h) (E) figures, plus two rows (Chainlink/Cosmos and
# Aptos/Aave) whose ranking order flipped and so are now swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.682.91"
$ws.Range("E2").Value = "  -3.56%  "

$ws.Range("D3").Value = "'1.743.80"
$ws.Range("E3").Value = "  -5.75%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'235.59"

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'0.4925"
$ws.Range("E7").Value = "  -7.85%  "

$ws.Range("D8").Value = "'41.47"
$ws.Range("E8").Value = "  -8.19%  "

$ws.Range("D9").Value = "'0.2540"
$ws.Range("E9").Value = "  -19.60%  "

$ws.Range("D10").Value = "'0.06025"
$ws.Range("E10").Value = "  -13.50%  "

$ws.Range("D11").Value = "'1.744.43"
$ws.Range("E11").Value = "  -5.79%  "

$ws.Range("D12").Value = "'0.06842"
$ws.Range("E12").Value = "  -12.58%  "

$ws.Range("D13").Value = "'14.80"
$ws.Range("E13").Value = "  -21.74%  "

$ws.Range("D14").Value = "'4.450"
$ws.Range("E14").Value = "  -11.90%  "

$ws.Range("D15").Value = "'76.78"
$ws.Range("E15").Value = "  -14.48%  "

$ws.Range("D16").Value = "'0.5673"
$ws.Range("E16").Value = "  -26.49%  "

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "'25.723.45"
$ws.Range("E19").Value = "  -3.48%  "

$ws.Range("D20").Value = "'11.28"
$ws.Range("E20").Value = "  -20.21%  "

$ws.Range("D21").Value = "'0.000006580"
$ws.Range("E21").Value = "  -17.53%  "

$ws.Range("D22").Value = "'1.968.61"
$ws.Range("E22").Value = "  -6.03%  "

$ws.Range("E23").Value = "  -13.92%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'5.039"
$ws.Range("E24").Value = "  -16.40%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'7.908"
$ws.Range("E25").Value = "  -15.55%  "

$ws.Range("D26").Value = "'137.30"
$ws.Range("E26").Value = "  -3.17%  "

$ws.Range("D27").Value = "'1.479"
$ws.Range("E27").Value = "  -12.89%  "

$ws.Range("D28").Value = "'1.826"
$ws.Range("E28").Value = "  -17.65%  "

$ws.Range("E29").Value = "  -14.18%  "

$ws.Range("D30").Value = "'101.94"
$ws.Range("E30").Value = "  -8.61%  "

$ws.Range("D31").Value = "'3.775"
$ws.Range("E31").Value = "  -12.35%  "

$ws.Range("D32").Value = "'0.07994"
$ws.Range("E32").Value = "  -8.91%  "

$ws.Range("D33").Value = "'3.407"
$ws.Range("E33").Value = "  -17.22%  "

$ws.Range("D34").Value = "'0.04387"
$ws.Range("E34").Value = "  -9.70%  "

$ws.Range("D35").Value = "'0.9992"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").Value = "'2.601"
$ws.Range("E36").Value = "  -9.66%  "

$ws.Range("D37").Value = "'0.9859"
$ws.Range("E37").Value = "  -13.55%  "

$ws.Range("D38").Value = "'0.6025"
$ws.Range("E38").Value = "  -18.28%  "

$ws.Range("D39").Value = "'2.679"
$ws.Range("E39").Value = "  -13.88%  "

$ws.Range("D40").Value = "'1.986"
$ws.Range("E40").Value = "  -15.91%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "'0.01509"
$ws.Range("E42").Value = "  -13.12%  "

$ws.Range("D43").Value = "'102.03"
$ws.Range("E43").Value = "  -6.36%  "

$ws.Range("D44").Value = "'0.7586"
$ws.Range("E44").Value = "  -16.82%  "

$ws.Range("D45").Value = "'5.173"
$ws.Range("E45").Value = "  -12.55%  "

$ws.Range("D46").Value = "'0.3740"
$ws.Range("E46").Value = "  -22.71%  "

$ws.Range("D47").Value = "'0.05226"
$ws.Range("E47").Value = "  -10.14%  "

$ws.Range("D48").Value = "'0.1065"
$ws.Range("E48").Value = "  -14.88%  "

$ws.Range("D49").Value = "'30.13"
$ws.Range("E49").Value = "  -14.08%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'52.32"
$ws.Range("E50").Value = "  -13.41%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'5.846"
$ws.Range("E51").Value = "  -24.20%  "
